# Clean up vaccine-name labels across every sheet of the workbook:
#   - collapse embedded line breaks (two-line cell text) into a single
#     line separated by a space
#   - drop the trailing footnote markers like " [1]" .. " [5]"
#
# This mirrors the author's own description ("fixed up the footnote /
# line-wrap text") and is applied uniformly to every cell that contains
# either a line break or a bracketed single-digit footnote, on every
# worksheet in the workbook.

$wb = $excel.ActiveWorkbook

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $ur = $ws.UsedRange

    $firstRow = $ur.Row
    $firstCol = $ur.Column
    $numRows = $ur.Rows.Count
    $numCols = $ur.Columns.Count
    $lastRow = $firstRow + $numRows - 1
    $lastCol = $firstCol + $numCols - 1

    for ($r = $firstRow; $r -le $lastRow; $r++) {
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $orig = $cell.Text

            if ($orig -ne $null -and ($orig -match "\[\d\]" -or $orig -match "`n")) {
                $updated = $orig -replace "`n", " "
                $updated = $updated -replace '\[\d\]', ''

                if ($updated -ne $orig) {
                    $cell.Value = $updated
                }
            }
        }
    }
}
